$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Mark the two completed tasks ("Estado" column) with an "x" to indicate
# the order ("pedido") process has been completed.
$ws.Range("C13").Value = "x"
$ws.Range("C14").Value = "x"

$ws.Range("C11").Select()
